# Target diff summary (public/template.xlsx):
#   1) workbook.xml: <sheet name="Magh" .../> -> <sheet name="Sheet" .../>
#   2) drawing1.xml: the embedded logo picture's anchor/extent shift by a
#      uniform, sub-pixel 360 EMU (0.01 mm / ~0.03pt) on all four numbers
#      (colOff, rowOff, ext cx, ext cy) - a classic LibreOffice-roundtrip
#      rounding artifact, not a deliberate resize/move (nothing in the
#      commit message relates to the logo). This headless COM surface
#      recomputes a shape's two-cell anchor from a fixed default row/col
#      grid rather than this sheet's real custom column widths, so poking
#      Shape.Width/Height/Left/Top here only trades a tiny, invisible
#      360 EMU drift for a much larger, visible corruption (e.g. the
#      anchor column flipping from col C to col B, or the picture's
#      width/height jumping by hundreds of thousands of EMU). Leaving the
#      picture untouched keeps it at the smallest possible, imperceptible
#      deviation from the target and avoids introducing a real regression.
#   3) sheet1.xml: selection activeCell/sqref I19 -> W22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet "Magh" -> "Sheet" (sheetId/position/relationship
#    stay the same - this is an in-place rename).
$ws.Name = "Sheet"

# 3) Move the selection/active cell to W22.
$ws.Range("W22").Select()
